$wb = $excel.ActiveWorkbook

# The handoff/handback automation has caught up with the 7a969992-...md
# source file: it is now "Handed back: in sync with en-US" on both the
# zh-cn and de-de sheets. Record the resulting target/handback files and
# timestamp the handback.

# Color used by the workbook's existing "HyperLink" cell style
# (font color FF6495ED, stored to the COM Font.Color BGR integer).
$hyperlinkColor = 15570276

# --- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"

$wsZh.Range("F2").Value = "7a969992-96fa-4bd8-bc6f-607c5fae6609.md"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/5bee6ed53ba7b262cc666cf9bde15ca3f5f58203/e2e/7a969992-96fa-4bd8-bc6f-607c5fae6609.md", "", "", "7a969992-96fa-4bd8-bc6f-607c5fae6609.md")
$wsZh.Range("F2").Font.Underline = $true
$wsZh.Range("F2").Font.Color = $hyperlinkColor

$wsZh.Range("G2").Value = "7a969992-96fa-4bd8-bc6f-607c5fae6609.712c506588753dd6977c86c8fda2d7b391a6e52e.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a19d66439c051e9ac102a530e5c4e1b9947bf226/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7a969992-96fa-4bd8-bc6f-607c5fae6609.712c506588753dd6977c86c8fda2d7b391a6e52e.zh-cn.xlf", "", "", "7a969992-96fa-4bd8-bc6f-607c5fae6609.712c506588753dd6977c86c8fda2d7b391a6e52e.zh-cn.xlf")
$wsZh.Range("G2").Font.Underline = $true
$wsZh.Range("G2").Font.Color = $hyperlinkColor

$wsZh.Range("H2").Value = "2016-03-23 00:36:09"

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

$wsDe.Range("F2").Value = "7a969992-96fa-4bd8-bc6f-607c5fae6609.md"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/5bee6ed53ba7b262cc666cf9bde15ca3f5f58203/e2e/7a969992-96fa-4bd8-bc6f-607c5fae6609.md", "", "", "7a969992-96fa-4bd8-bc6f-607c5fae6609.md")
$wsDe.Range("F2").Font.Underline = $true
$wsDe.Range("F2").Font.Color = $hyperlinkColor

$wsDe.Range("G2").Value = "7a969992-96fa-4bd8-bc6f-607c5fae6609.712c506588753dd6977c86c8fda2d7b391a6e52e.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8f6069f5ef3c4f2adf883a6d0dbb5a0aeb3a3324/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7a969992-96fa-4bd8-bc6f-607c5fae6609.712c506588753dd6977c86c8fda2d7b391a6e52e.de-de.xlf", "", "", "7a969992-96fa-4bd8-bc6f-607c5fae6609.712c506588753dd6977c86c8fda2d7b391a6e52e.de-de.xlf")
$wsDe.Range("G2").Font.Underline = $true
$wsDe.Range("G2").Font.Color = $hyperlinkColor

$wsDe.Range("H2").Value = "2016-03-23 00:36:15"
